$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A day of data (August 5, 2025) was missing from the daily revenue table.
# Insert a new row at row 6 (pushing the existing July/June/May blocks down
# by one row) and fill it with the missing day's figures.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 21047.97
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 2025
$ws.Range("E6").Value = "08/2025"
